$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("potential")

# New row 40: BECCS / 2.5 / Hanssen 2020
# Values mirror the existing text-typed numbers already present in this sheet
# (e.g. B8/B13 hold "2.5" as text), so a leading apostrophe keeps them text
# instead of auto-converting to numeric cells. Applying it uniformly to the
# whole new row keeps A/B/C on one consistent cell style.
$ws.Range("A40").Value = "'BECCS"
$ws.Range("B40").Value = "'2.5"
$ws.Range("C40").Value = "'Hanssen 2020"

# New row 41: BECCS / 2.8 / Pour 2018
$ws.Range("A41").Value = "'BECCS"
$ws.Range("B41").Value = "'2.8"
$ws.Range("C41").Value = "'Pour 2018"

# Match the final selection left behind by the edit session.
$ws.Range("F34").Select()
